$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 1 (SL. No "1.") explanation cell: collapse the "DropBox" spell-check
# split runs into a single plain run with no w:proofErr markers. Replacing a
# range that spans across all three original runs (and the proofErr tags
# sitting between them) merges them into one run and drops the stray
# proofErr elements.
$explanationCell = $t.Cell(2, 3).Range
$find = $explanationCell.Find
$find.Execute("the DropBox API", $false, $false, $false, $false, $false, `
    $true, 1, $false, "the DropBox API", 2) | Out-Null

# --- Row 2 (previously empty placeholder row): fill in SL. No, Error Code,
# and Explanation text.
$t.Cell(3, 1).Range.Text = "2."
$t.Cell(3, 2).Range.Text = "101"
$t.Cell(3, 3).Range.Text = "Could not read the Excel File to get the counter"
